$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 14-17: summary labels + aggregate formulas
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

foreach ($addr in @("B14", "B15", "B16", "B17")) {
    $r = $ws.Range($addr)
    $r.Font.Bold = $true
    $r.Font.Size = 12
    $r.VerticalAlignment = -4108
}

$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Row 12: average of the k column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("J12").Select() | Out-Null
